# Fruta / hortaliza, semanal
# Insert a new weekly price record as row 83 in the "Poroto verde" sheet,
# pushing the existing rows 83-149 down to 84-150.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the current row 83 (shifts 83..149 -> 84..150,
# carrying the existing row formatting, e.g. the date style on column D).
$ws.Rows.Item(83).Insert()

# Populate the newly inserted row 83 with the new record.
$ws.Cells.Item(83, 1).Value = 4
$ws.Cells.Item(83, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(83, 3).Value = "Los Lagos"
$ws.Cells.Item(83, 4).Value = 45090
$ws.Cells.Item(83, 5).Value = 10
$ws.Cells.Item(83, 6).Value = 100112031
$ws.Cells.Item(83, 7).Value = "Poroto verde"
$ws.Cells.Item(83, 8).Value = "Magnum"
$ws.Cells.Item(83, 9).Value = "Primera"
$ws.Cells.Item(83, 10).Value = 45
$ws.Cells.Item(83, 11).Value = 27000
$ws.Cells.Item(83, 12).Value = 27000
$ws.Cells.Item(83, 13).Value = 27000
$ws.Cells.Item(83, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(83, 15).Value = "Perú"
$ws.Cells.Item(83, 16).Value = 1080
$ws.Cells.Item(83, 17).Value = 25
$ws.Cells.Item(83, 18).Value = "Hortaliza"
